$d = $word.ActiveDocument

function Replace-Exact($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Text = $find
    $range.Find.Replacement.Text = $replace
    $range.Find.Forward = $true
    $range.Find.Wrap = 1
    $range.Find.Format = $false
    $range.Find.MatchCase = $true
    $range.Find.MatchWholeWord = $false
    $range.Find.MatchWildcards = $false
    $range.Find.MatchSoundsLike = $false
    $range.Find.MatchAllWordForms = $false
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# RA number cleared
Replace-Exact " 000109085235 - 6 " "  "

# Case-sensitive placeholder replacements (document body + header)
Replace-Exact "QWER" "TRE"
Replace-Exact "QWR" "TERE"
Replace-Exact "Qwer" "Tre"
Replace-Exact "qwer" "tre"
